$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1.31

$ws.Range("D3").Value = 1.42
$ws.Range("F3").Value = 1.23

$ws.Range("B4").Value = 1.4
$ws.Range("C4").Value = 1.42
$ws.Range("D4").Value = 1.34
$ws.Range("F4").Value = 1.07

$ws.Range("C5").Value = 1.37
$ws.Range("E5").Value = 1.22
$ws.Range("F5").Value = 1.06
$ws.Range("G5").Value = 0.65

$ws.Range("C6").Value = 1.46
$ws.Range("D6").Value = 1.53
$ws.Range("E6").Value = 1.3
$ws.Range("G6").Value = 1.04

$ws.Range("E7").Value = 1.96
$ws.Range("F7").Value = 1.5
